$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; this shifts existing rows 41-78 down to 42-79,
# carrying all their data/formatting with them.
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new data record.
$ws.Range("A41").Value = 9
$ws.Range("B41").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44778
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100102
$ws.Range("H41").Value = "Cítricos"
$ws.Range("I41").Value = 100102006
$ws.Range("J41").Value = "Pomelo"
$ws.Range("K41").Value = "Start Ruby"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 400
$ws.Range("N41").Value = 8000
$ws.Range("O41").Value = 8000
$ws.Range("P41").Value = 8000
$ws.Range("Q41").Value = "$/caja 14 kilos"
$ws.Range("R41").Value = "Región Metropolitana"
$ws.Range("S41").Value = 571
$ws.Range("T41").Value = 14
